$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the PAYOUT% (I) and PAYOUT (J) columns entirely - shifts nothing right of J,
# collapses dimension from J to H.
$ws.Range("I1:J12").Delete()

# Update vintage labels for the PL Self block (rows 10-12 shift down one vintage,
# V1 data stays, old V3 row removed upstream so V4->V3, V5->V4, V6->V5)
$ws.Range("B10").Value = "V3"
$ws.Range("B11").Value = "V4"
$ws.Range("B12").Value = "V5"

# Row 2: PL Sal / V1
$ws.Range("C2").Value = 9170193.810000001
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 161172
$ws.Range("H2").Value = 1.76

# Row 3: PL Sal / V2
$ws.Range("C3").Value = 17749604.12
$ws.Range("D3").Value = 52
$ws.Range("E3").Value = 47
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 181585
$ws.Range("H3").Value = 1.02

# Row 4: PL Sal / V3
$ws.Range("C4").Value = 1448782.22
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5000
$ws.Range("H4").Value = 0.35

# Row 5: PL Sal / V4
$ws.Range("C5").Value = 10197523.39
$ws.Range("D5").Value = 31
$ws.Range("E5").Value = 29
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 30000
$ws.Range("H5").Value = 0.29

# Row 6: PL Sal / V5
$ws.Range("C6").Value = 13244323.61
$ws.Range("D6").Value = 41
$ws.Range("E6").Value = 40
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15000
$ws.Range("H6").Value = 0.11

# Row 7: PL Sal / V6
$ws.Range("C7").Value = 1828497.2
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 25000
$ws.Range("H7").Value = 1.37

# Row 8: PL Self / V1
$ws.Range("C8").Value = 2805254.55
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15000
$ws.Range("H8").Value = 0.53

# Row 9: PL Self / V2
$ws.Range("C9").Value = 7479921.47
$ws.Range("D9").Value = 15
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 30000
$ws.Range("H9").Value = 0.4

# Row 10: PL Self / V3 (was V4)
$ws.Range("C10").Value = 550961.3100000001
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

# Row 11: PL Self / V4 (was V5)
$ws.Range("C11").Value = 4588615.25
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10000
$ws.Range("H11").Value = 0.22

# Row 12: PL Self / V5 (was V6)
$ws.Range("C12").Value = 4249598.26
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

# Row 13 (new): PL Self / V6
$ws.Range("A13").Value = "PL Self"
$ws.Range("B13").Value = "V6"
$ws.Range("C13").Value = 1818015.45
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 32000
$ws.Range("H13").Value = 1.76
